$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "0105054971186"
$ws.Range("C6").Value = "ID Used"
$ws.Range("A7").Value = "98050503751810"
$ws.Range("C7").Value = "ID Used"
$ws.Range("A8").Value = "9505054540082"
$ws.Range("C8").Value = "ID Used"
$ws.Range("A9").Value = "0105052723183"
$ws.Range("C9").Value = "ID Used"
$ws.Range("A10").Value = "9805053170083"
$ws.Range("C10").Value = "ID Used"
$ws.Range("A11").Value = "9505054380182"
$ws.Range("C11").Value = "ID Used"
$ws.Range("A12").Value = "9405053071082"
$ws.Range("C12").Value = "ID Used"
$ws.Range("A13").Value = "0105051312087"
$ws.Range("C13").Value = "ID Used"
$ws.Range("A14").Value = "9805051123084"
$ws.Range("C14").Value = "ID Used"
$ws.Range("A15").Value = "9505052045183"
$ws.Range("C15").Value = "ID Used"
$ws.Range("A16").Value = "9105054038185"
$ws.Range("C16").Value = "ID Used"
$ws.Range("A17").Value = "0105054842189"
$ws.Range("C17").Value = "ID Used"
$ws.Range("A18").Value = "9805050603086"
$ws.Range("C18").Value = "ID Used"
$ws.Range("A19").Value = "9505050077188"
$ws.Range("C19").Value = "ID Used"
$ws.Range("A20").Value = "9405052188085"
$ws.Range("C20").Value = "ID Used"
$ws.Range("A21").Value = "0105053886187"
$ws.Range("C21").Value = "ID Used"
$ws.Range("A22").Value = "9805050383184"
$ws.Range("C22").Value = "ID Used"
$ws.Range("A23").Value = "9505050874089"
$ws.Range("C23").Value = "ID Used"
$ws.Range("A24").Value = "0105054928087"
$ws.Range("C24").Value = "ID Used"
$ws.Range("A25").Value = "9805052024083"
$ws.Range("C25").Value = "ID Used"
$ws.Range("A26").Value = "9505053421086"
$ws.Range("C26").Value = "ID Used"
$ws.Range("A27").Value = "0105053024086"
$ws.Range("C27").Value = "ID Used"
$ws.Range("A28").Value = "9805054865087"
$ws.Range("C28").Value = "ID Used"
$ws.Range("A29").Value = "9505050468189"
$ws.Range("C29").Value = "ID Used"
$ws.Range("A30").Value = "0105054828188"
$ws.Range("C30").Value = "ID Used"
$ws.Range("A31").Value = "9805052011189"
$ws.Range("C31").Value = "ID Used"
$ws.Range("A32").Value = "9505051407087"
$ws.Range("C32").Value = "ID Used"

$null = $ws.Range("U11").Select()
